# Update the "南宁·9.7国乙同人ONLY" event (now cancelled) on both the
# "展览" (exhibitions) sheet and the "全部类型" (all types) sheet.
#
# Changes per affected sheet:
#   C2: "南宁·9.7国乙同人ONLY" -> "南宁·9.7国乙同人ONLY（取消）"
#   F2: 43   -> 44
#   G2: 109  -> "不可售"   (becomes a text value, no longer numeric)
#   F4: 2131 -> 2160
#   F5: 173  -> 182
# Plus, on each sheet, the last data row's F column (F6 on 展览, F7 on 全部类型):
#   361 -> 362

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Cancelled event name
    $ws.Range("C2").Value = "南宁·9.7国乙同人ONLY（取消）"

    # Updated interest count
    $ws.Range("F2").Value = 44

    # Minimum price is no longer sellable -> text value
    $ws.Range("G2").Value = "不可售"

    # Interest counts bumped on a couple other rows
    $ws.Range("F4").Value = 2160
    $ws.Range("F5").Value = 182

    # Last data row differs per sheet (6 on 展览, 7 on 全部类型)
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $ws.Cells.Item($lastRow, 6).Value = 362
}
